# Apply cryptos list update (commit: "Updated cryptos list on Sun Oct 20 04:33:46 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: a leading apostrophe forces Excel to store the value
# as literal text (matching the original inlineStr cells) instead of letting
# it auto-coerce numeric-looking strings (e.g. "1.00", "596.39") into numbers.
# Resetting the cell Style back to "Normal" afterwards drops the quotePrefix
# style Excel would otherwise stamp on the cell, keeping it unstyled like the source.

# --- Rows 2-41: refresh Price (D) and Volume(1h) (E) figures ---
$ws.Range("D2").Value = "'68.286.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "'2.645.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'596.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "'159.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.22%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.543"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("D9").Value = "'0.143"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("D11").Value = "'5.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "'0.351"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").Value = "'27.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "'3.130.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "'0.0000187"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").Value = "'68.200.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "'2.618.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "'11.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "'360.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "'7.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").Value = "'4.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "'4.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("D23").Value = "'2.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").Value = "'75.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'9.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'2.775.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "'0.0000104"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'563.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("D31").Value = "'7.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("D32").Value = "'1.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("D33").Value = "'1.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("D36").Value = "'1.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("D37").Value = "'19.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("D38").Value = "'158.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "'0.370"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").Value = "'1.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("D41").Value = "'5.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.03%  "

# --- Rows 42-51: WhiteBITCoin dropped from the list; remaining coins shift up
#     one row and ARBITRUM is newly appended as row 51 ---
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("B43").Value = "BabyDogeCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D43").Value = "'0.0₆0323"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.27%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'157.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'3.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'21.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("B48").Value = "Optimism"
$ws.Range("C48").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D48").Value = "'1.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0774"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.612"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'0.566"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.18%  "
